$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A. This shifts every existing cell,
# the merged ranges, and the column width definitions one column to the
# right (B->C, C->D, ... J->K), matching the target layout.
$ws.Columns("A:A").Insert()

# The rounded-rectangle shape keeps the exact same cell-relative offset it
# had from (old) column A, but now that a brand-new column A sits in front
# of it, its absolute position on the sheet must shift right by that new
# column's width so it keeps starting at the same spot inside (new) column
# B and ending at the same spot inside (new) column L. The figures below
# are the sheet's original anchor (from-col 0 @ 90488 EMU offset, to-col 10
# @ 4763 EMU offset) re-expressed in points at full precision so the
# re-anchored shape lands on the exact same pixel-accurate column offsets,
# just one column further right.
$shp = $ws.Shapes.Item(1)
$shp.Left = 65.56253937007874
$shp.Width = 654.7001953125

# Match the workbook's recorded selection after the edit.
$ws.Range("M3").Select()
